$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Record the new state/region name used for this copy of the model
$ws.Range("B1").Value = "Colorado"

# Update the "last updated" date shown in cell C1 (11/9/2021)
$ws.Range("C1").Value = 44509
